# Generate Report for Handoff
# Adds a new tracked file (8ee5f43b-594c-4ef3-86de-e29f4246037d.md) as row 3
# to the Overview, zh-cn and de-de sheets, mirroring the existing row for
# 8311c0ce-0600-4c2c-b051-9dd523a223c9.md.

$wb = $excel.ActiveWorkbook

$newGuid = "8ee5f43b-594c-4ef3-86de-e29f4246037d"
$newHash = "932e5a8e2614e9da113243bc511486492f4eae5e"
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1093fcfccc3703e1da97d39cc820c0a4700545a2/e2e/$newGuid.md"
$hyperlinkDisplay = "e2e\$newGuid.md"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) - table3 / displayName "Overview"
# Columns: File Name | Path And Name | Extension | Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()
$rngOverview = $rowOverview.Range

$rngOverview.Item(1, 1).Value = "$newGuid.md"
$rngOverview.Item(1, 3).Value = ".md"
$rngOverview.Item(1, 4).Value = ""
$rngOverview.Item(1, 5).Value = "Ready for handoff"
$rngOverview.Item(1, 6).Value = "Ready for handoff"
$rngOverview.Item(1, 7).Value = "2016-09-01 16:47:50"
$rngOverview.Item(1, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($rngOverview.Item(1, 2), $hyperlinkUrl, "", "", $hyperlinkDisplay)

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - table1
# Columns: Source File Name | File Extension | Status | Source Path | Priority |
#          Content Duplicate | Latest Handoff File | Latest Handoff Datetime |
#          Latest Target File | Latest Handback File | Latest Handback DateTime |
#          Reference Tokens | To be localized | Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()
$rngZhCn = $rowZhCn.Range

$rngZhCn.Item(1, 2).Value = ".md"
$rngZhCn.Item(1, 3).Value = "Ready for handoff"
$rngZhCn.Item(1, 4).Value = "e2e"
$rngZhCn.Item(1, 5).Value = "ht"
$rngZhCn.Item(1, 6).Value = "False"
$rngZhCn.Item(1, 7).Value = "$newGuid.$newHash.zh-cn.xlf"
$rngZhCn.Item(1, 8).Value = "2016-09-01 16:47:45"
$rngZhCn.Item(1, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$rngZhCn.Item(1, 9).Value = ""
$rngZhCn.Item(1, 10).Value = ""
$rngZhCn.Item(1, 11).Value = "0001-01-01 00:00:00"
$rngZhCn.Item(1, 12).Value = ""
$rngZhCn.Item(1, 13).Value = "True"
$rngZhCn.Item(1, 14).Value = ""
$rngZhCn.Item(1, 15).Value = "False"
$rngZhCn.Item(1, 16).Value = ""

$wsZhCn.Hyperlinks.Add($rngZhCn.Item(1, 1), $hyperlinkUrl, "", "", "$newGuid.md")

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) - table2 (same column layout as zh-cn)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()
$rngDeDe = $rowDeDe.Range

$rngDeDe.Item(1, 2).Value = ".md"
$rngDeDe.Item(1, 3).Value = "Ready for handoff"
$rngDeDe.Item(1, 4).Value = "e2e"
$rngDeDe.Item(1, 5).Value = "ht"
$rngDeDe.Item(1, 6).Value = "False"
$rngDeDe.Item(1, 7).Value = "$newGuid.$newHash.de-de.xlf"
$rngDeDe.Item(1, 8).Value = "2016-09-01 16:47:50"
$rngDeDe.Item(1, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$rngDeDe.Item(1, 9).Value = ""
$rngDeDe.Item(1, 10).Value = ""
$rngDeDe.Item(1, 11).Value = "0001-01-01 00:00:00"
$rngDeDe.Item(1, 12).Value = ""
$rngDeDe.Item(1, 13).Value = "True"
$rngDeDe.Item(1, 14).Value = ""
$rngDeDe.Item(1, 15).Value = "False"
$rngDeDe.Item(1, 16).Value = ""

$wsDeDe.Hyperlinks.Add($rngDeDe.Item(1, 1), $hyperlinkUrl, "", "", "$newGuid.md")
